# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Mon Apr  8 14:34:47 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be stored as text (not auto-converted to a number)
    # by temporarily applying a text number format, then restore the default style
    # so the cell ends up with no explicit style, matching the original workbook.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "71.563.38"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "3.634.86"
$ws.Range("E3").Value = "  +6.99%  "
$ws.Range("E4").Value = "  +0.20%  "
Set-TextValue $ws.Range("D5") "589.62"
$ws.Range("E5").Value = "  +0.30%  "
Set-TextValue $ws.Range("D6") "180.21"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "3.625.94"
$ws.Range("E7").Value = "  +7.02%  "
Set-TextValue $ws.Range("D8") "0.611"
$ws.Range("E8").Value = "  +2.05%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -0.03%  "
Set-TextValue $ws.Range("D11") "0.605"
$ws.Range("E11").Value = "  +2.03%  "
Set-TextValue $ws.Range("D12") "49.69"
$ws.Range("E12").Value = "  +2.56%  "
Set-TextValue $ws.Range("D13") "0.0000285"
$ws.Range("E13").Value = "  -0.61%  "
Set-TextValue $ws.Range("D14") "684.56"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "4.217.41"
$ws.Range("E15").Value = "  +6.84%  "
Set-TextValue $ws.Range("D16") "8.98"
$ws.Range("E16").Value = "  +3.29%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.659.47"
$ws.Range("E17").Value = "  +7.71%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "71.765.96"
$ws.Range("E18").Value = "  +3.52%  "
$ws.Range("E19").Value = "  +1.82%  "
Set-TextValue $ws.Range("D20") "18.29"
$ws.Range("E20").Value = "  +3.20%  "
Set-TextValue $ws.Range("D21") "11.59"
$ws.Range("E21").Value = "  +2.37%  "
Set-TextValue $ws.Range("D22") "0.934"
$ws.Range("E22").Value = "  +2.41%  "
Set-TextValue $ws.Range("D23") "5.85"
$ws.Range("E23").Value = "  +9.26%  "
Set-TextValue $ws.Range("D24") "17.76"
$ws.Range("E24").Value = "  +2.61%  "
Set-TextValue $ws.Range("D25") "103.20"
$ws.Range("E25").Value = "  -0.14%  "
Set-TextValue $ws.Range("D26") "4.01"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("E27").Value = "  +3.88%  "
Set-TextValue $ws.Range("D28") "9.97"
Set-TextValue $ws.Range("D29") "34.99"
$ws.Range("E29").Value = "  +3.24%  "
Set-TextValue $ws.Range("D30") "9.18"
$ws.Range("E30").Value = "  +4.26%  "
Set-TextValue $ws.Range("D31") "7.24"
$ws.Range("E31").Value = "  +4.07%  "
Set-TextValue $ws.Range("D32") "4.19"
$ws.Range("E32").Value = "  +15.66%  "
Set-TextValue $ws.Range("D33") "578.71"
$ws.Range("E33").Value = "  +3.90%  "
Set-TextValue $ws.Range("D34") "11.32"
$ws.Range("E34").Value = "  +1.40%  "
Set-TextValue $ws.Range("D35") "0.109"
$ws.Range("E35").Value = "  +2.49%  "
Set-TextValue $ws.Range("D36") "59.43"
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "3.671.72"
$ws.Range("E38").Value = "  +0.46%  "
Set-TextValue $ws.Range("D39") "0.142"
$ws.Range("E39").Value = "  +0.38%  "
Set-TextValue $ws.Range("D40") "35.70"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").Value = "0.0₃0760"
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D42") "0.0475"
$ws.Range("E42").Value = "  +11.09%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D43") "3.40"
$ws.Range("E43").Value = "  +4.52%  "
Set-TextValue $ws.Range("D44") "2.75"
$ws.Range("E44").Value = "  +3.11%  "
Set-TextValue $ws.Range("D45") "0.347"
$ws.Range("E45").Value = "  +2.48%  "
Set-TextValue $ws.Range("D46") "3.38"
$ws.Range("E46").Value = "  +1.11%  "
Set-TextValue $ws.Range("D47") "2.82"
$ws.Range("E47").Value = "  +5.52%  "
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("E49").Value = "  +3.15%  "
$ws.Range("E50").Value = "  +0.09%  "
Set-TextValue $ws.Range("D51") "131.90"
$ws.Range("E51").Value = "  +0.41%  "
